# Sprint Twelve Team Backlog update:
# Insert 5 new SYSTICK_* backlog rows above the existing "Main Flow App"
# row, pushing the trailing rows (Main Flow App / Testing All Apis /
# Documentation) down, and touch up a couple of their Assigned/time cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing "Main Flow App" row (row 14) is where the new tasks need to
# land; push it (and the two rows below it) down by 5 rows first.
$ws.Rows("14:18").Insert()

# --- New rows 14-18: SYSTICK driver tasks -------------------------------
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = "SYSTICK_init"
$ws.Range("C14").Value = "Sherif Ashraf"
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = "1 hour"
$ws.Range("F14").Value = "1 hour"

$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "SYSTICK_interruptEnable"
$ws.Range("C15").Value = "Sherif Ashraf"
$ws.Range("D15").Value = 100
$ws.Range("E15").Value = "30 minutes"
$ws.Range("F15").Value = "15 minutes"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "SYSTICK_interruptDisable"
$ws.Range("C16").Value = "Momen Hassan"
$ws.Range("D16").Value = 100
$ws.Range("E16").Value = "30 minutes"
$ws.Range("F16").Value = "15 minutes"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "SYSTICK_setDelayInMs"
$ws.Range("C17").Value = "Momen Hassan"
$ws.Range("D17").Value = 100
$ws.Range("E17").Value = "30 minutes"
$ws.Range("F17").Value = "15 minutes"

$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "SYSTICK_handler"
$ws.Range("C18").Value = "Momen Hassan"
$ws.Range("D18").Value = 100
$ws.Range("E18").Value = "15 minutes"
$ws.Range("F18").Value = "15 minutes"

# --- Rows 19-21: renumber + touch up the pushed-down rows ---------------
$ws.Range("A19").Value = 18
$ws.Range("C19").Value = "Sherif Ashraf"

$ws.Range("A20").Value = 19

$ws.Range("A21").Value = 20
$ws.Range("C21").Value = "All"

# --- Conditional formatting: extend the data-bar range by the 5 new rows
$fc = $ws.Range("D2:D51").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("D2:D56"))
